$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "time_taken" column, matching style of existing headers (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Timestamps captured for each data row (rows 2-38), to be written as plain text
$timestamps = @("2021-10-05 10:50:32.953744","2021-10-05 10:50:32.953755","2021-10-05 10:50:32.953759","2021-10-05 10:50:32.953761","2021-10-05 10:50:32.953764","2021-10-05 10:50:32.953767","2021-10-05 10:50:32.953770","2021-10-05 10:50:32.953773","2021-10-05 10:50:32.953776","2021-10-05 10:50:32.953778","2021-10-05 10:50:32.953781","2021-10-05 10:50:32.953783","2021-10-05 10:50:32.953786","2021-10-05 10:50:32.953789","2021-10-05 10:50:32.953791","2021-10-05 10:50:32.953794","2021-10-05 10:50:32.953797","2021-10-05 10:50:32.953800","2021-10-05 10:50:32.953802","2021-10-05 10:50:32.953805","2021-10-05 10:50:32.953808","2021-10-05 10:50:32.953810","2021-10-05 10:50:32.953813","2021-10-05 10:50:32.953815","2021-10-05 10:50:32.953818","2021-10-05 10:50:32.953821","2021-10-05 10:50:32.953823","2021-10-05 10:50:32.953826","2021-10-05 10:50:32.953828","2021-10-05 10:50:32.953831","2021-10-05 10:50:32.953834","2021-10-05 10:50:32.953837","2021-10-05 10:50:32.953840","2021-10-05 10:50:32.953843","2021-10-05 10:50:32.953845","2021-10-05 10:50:32.953848","2021-10-05 10:50:32.953850")

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
